$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "42.487.98"
$ws.Range("E2").Value = "  -1.17%  "

# Row 3
$ws.Range("D3").Value = "2.219.97"
$ws.Range("E3").Value = "  -0.82%  "

# Row 5
$ws.Range("D5").Value = "'111.16"
$ws.Range("E5").Value = "  -3.98%  "

# Row 6
$ws.Range("D6").Value = "'289.66"
$ws.Range("E6").Value = "  +9.15%  "

# Row 7
$ws.Range("E7").Value = "  -0.61%  "

# Row 8
$ws.Range("E8").Value = "  -0.50%  "

# Row 9
$ws.Range("E9").Value = "  -1.55%  "

# Row 10
$ws.Range("E10").Value = "  -6.26%  "

# Row 11
$ws.Range("D11").Value = "'0.0910"
$ws.Range("E11").Value = "  -1.70%  "

# Row 12
$ws.Range("D12").Value = "'54.31"
$ws.Range("E12").Value = "  +0.93%  "

# Row 13
$ws.Range("D13").Value = "'8.59"
$ws.Range("E13").Value = "  -5.92%  "

# Row 14
$ws.Range("D14").Value = "'1.01"
$ws.Range("E14").Value = "  +14.08%  "

# Row 15
$ws.Range("E15").Value = "  -1.84%  "

# Row 16
$ws.Range("D16").Value = "'14.84"
$ws.Range("E16").Value = "  -3.20%  "

# Row 17
$ws.Range("D17").Value = "2.553.85"
$ws.Range("E17").Value = "  -0.98%  "

# Row 18
$ws.Range("D18").Value = "2.229.41"
$ws.Range("E18").Value = "  -1.14%  "

# Row 19
$ws.Range("D19").Value = "42.344.00"
$ws.Range("E19").Value = "  -1.56%  "

# Row 20
$ws.Range("B20").Value = "ShibaInu"
$ws.Range("C20").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D20").Value = "'0.0000105"
$ws.Range("E20").Value = "  -1.73%  "

# Row 21
$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D21").Value = "'7.11"
$ws.Range("E21").Value = "  +6.07%  "

# Row 22
$ws.Range("D22").Value = "'73.25"
$ws.Range("E22").Value = "  +2.32%  "

# Row 23
$ws.Range("D23").Value = "'3.31"
$ws.Range("E23").Value = "  +14.68%  "

# Row 24
$ws.Range("D24").Value = "'2.39"
$ws.Range("E24").Value = "  +0.74%  "

# Row 25
$ws.Range("D25").Value = "'234.86"
$ws.Range("E25").Value = "  +1.60%  "

# Row 26
$ws.Range("D26").Value = "'8.87"
$ws.Range("E26").Value = "  -6.53%  "

# Row 27
$ws.Range("D27").Value = "'1.00"
$ws.Range("E27").Value = "  -1.86%  "

# Row 28
$ws.Range("E28").Value = "  -5.62%  "

# Row 29
$ws.Range("E29").Value = "  -1.69%  "

# Row 30
$ws.Range("B30").Value = "InjectiveProtocol"
$ws.Range("C30").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D30").Value = "'37.38"
$ws.Range("E30").Value = "  -8.99%  "

# Row 31
$ws.Range("B31").Value = "Monero"
$ws.Range("C31").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D31").Value = "'173.23"
$ws.Range("E31").Value = "  +0.38%  "

# Row 32
$ws.Range("D32").Value = "'3.12"
$ws.Range("E32").Value = "  -5.49%  "

# Row 33
$ws.Range("D33").Value = "'21.25"
$ws.Range("E33").Value = "  +0.45%  "

# Row 34
$ws.Range("D34").Value = "'0.0874"
$ws.Range("E34").Value = "  -2.40%  "

# Row 35
$ws.Range("D35").Value = "'5.58"
$ws.Range("E35").Value = "  +0.04%  "

# Row 36
$ws.Range("E36").Value = "  +7.10%  "

# Row 37
$ws.Range("D37").Value = "'0.125"
$ws.Range("E37").Value = "  -1.64%  "

# Row 38
$ws.Range("E38").Value = "  -3.55%  "

# Row 39
$ws.Range("D39").Value = "'0.0377"
$ws.Range("E39").Value = "  +1.67%  "

# Row 40
$ws.Range("E40").Value = "  -1.55%  "

# Row 41
$ws.Range("E41").Value = "  -4.12%  "

# Row 42
$ws.Range("D42").Value = "'71.48"
$ws.Range("E42").Value = "  +0.70%  "

# Row 43
$ws.Range("E43").Value = "  -2.55%  "

# Row 44
$ws.Range("E44").Value = "  -0.25%  "

# Row 45
$ws.Range("D45").Value = "'12.28"
$ws.Range("E45").Value = "  -7.66%  "

# Row 46
$ws.Range("D46").Value = "'1.30"
$ws.Range("E46").Value = "  -2.46%  "

# Row 47
$ws.Range("D47").Value = "'5.32"
$ws.Range("E47").Value = "  -5.28%  "

# Row 48
$ws.Range("E48").Value = "  +2.46%  "

# Row 49
$ws.Range("D49").Value = "'1.64"
$ws.Range("E49").Value = "  +4.22%  "

# Row 50
$ws.Range("B50").Value = "Aave"
$ws.Range("C50").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D50").Value = "'101.03"
$ws.Range("E50").Value = "  +1.04%  "

# Row 51
$ws.Range("B51").Value = "FraxShare"
$ws.Range("C51").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D51").Value = "'8.39"
$ws.Range("E51").Value = "  -0.12%  "
